# Applies the cryptos.xlsx price/volume/name/link updates described in the diff.
# Every target cell is stored as literal text in the source file (inlineStr), so we
# force text with a leading apostrophe via .Formula (Excel text-prefix convention)
# rather than .Value, which would otherwise coerce numeric-looking strings (e.g.
# "1.00" -> 1, "73.40" -> 73.4) and strip significant trailing zeros / formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $ws.Cells.Item($row, $col).Formula = "'" + $text
}

Set-TextCell 2 4 "37.009.16"
Set-TextCell 2 5 "  +1.02%  "
Set-TextCell 3 4 "2.064.31"
Set-TextCell 3 5 "  -1.59%  "
Set-TextCell 4 5 "  +0.04%  "
Set-TextCell 5 4 "249.72"
Set-TextCell 5 5 "  -1.32%  "
Set-TextCell 6 4 "0.674"
Set-TextCell 6 5 "  +2.06%  "
Set-TextCell 7 5 "  -0.06%  "
Set-TextCell 8 4 "55.24"
Set-TextCell 8 5 "  +11.67%  "
Set-TextCell 9 4 "61.09"
Set-TextCell 9 5 "  +1.33%  "
Set-TextCell 10 5 "  +1.55%  "
Set-TextCell 11 5 "  +7.33%  "
Set-TextCell 12 5 "  +5.81%  "
Set-TextCell 13 4 "15.06"
Set-TextCell 13 5 "  +1.79%  "
Set-TextCell 14 4 "2.362.91"
Set-TextCell 14 5 "  -1.73%  "
Set-TextCell 15 4 "0.818"
Set-TextCell 15 5 "  -2.23%  "
Set-TextCell 17 4 "2.056.73"
Set-TextCell 17 5 "  -1.86%  "
Set-TextCell 18 4 "37.002.23"
Set-TextCell 18 5 "  +1.05%  "
Set-TextCell 19 4 "0.0₃0943"
Set-TextCell 19 5 "  +12.72%  "
Set-TextCell 20 4 "73.40"
Set-TextCell 20 5 "  +0.30%  "
Set-TextCell 21 4 "14.21"
Set-TextCell 21 5 "  +6.61%  "
Set-TextCell 22 4 "5.41"
Set-TextCell 22 5 "  +2.20%  "
Set-TextCell 23 4 "237.67"
Set-TextCell 23 5 "  -1.16%  "
Set-TextCell 24 5 "  +0.00%  "
Set-TextCell 25 4 "2.44"
Set-TextCell 25 5 "  -3.97%  "
Set-TextCell 26 4 "170.94"
Set-TextCell 26 5 "  -0.33%  "
Set-TextCell 27 4 "9.13"
Set-TextCell 27 5 "  -1.44%  "
Set-TextCell 28 4 "20.16"
Set-TextCell 28 5 "  -5.07%  "
Set-TextCell 29 5 "  +0.65%  "
Set-TextCell 30 5 "  +1.86%  "
Set-TextCell 31 5 "  +2.36%  "
Set-TextCell 33 4 "0.0630"
Set-TextCell 33 5 "  +1.54%  "
Set-TextCell 34 4 "4.39"
Set-TextCell 34 5 "  +7.15%  "
Set-TextCell 35 4 "0.0893"
Set-TextCell 35 5 "  -0.56%  "
Set-TextCell 36 4 "1.00"
Set-TextCell 36 5 "  -0.06%  "
Set-TextCell 37 4 "2.28"
Set-TextCell 37 5 "  -6.01%  "
Set-TextCell 38 5 "  -4.71%  "
Set-TextCell 39 4 "1.35"
Set-TextCell 39 5 "  -0.60%  "
Set-TextCell 40 5 "  +24.18%  "
Set-TextCell 41 2 "InjectiveProtocol"
Set-TextCell 41 3 "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell 41 4 "17.78"
Set-TextCell 41 5 "  +7.36%  "
Set-TextCell 42 2 "VeChain"
Set-TextCell 42 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell 42 4 "0.0226"
Set-TextCell 42 5 "  +0.13%  "
Set-TextCell 43 5 "  -2.32%  "
Set-TextCell 44 4 "96.78"
Set-TextCell 44 5 "  -1.38%  "
Set-TextCell 46 4 "4.12"
Set-TextCell 46 5 "  +36.90%  "
Set-TextCell 47 4 "13.90"
Set-TextCell 47 5 "  -50.95%  "
Set-TextCell 48 5 "  +7.58%  "
Set-TextCell 49 2 "THORChain"
Set-TextCell 49 3 "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell 49 4 "4.24"
Set-TextCell 49 5 "  +9.50%  "
Set-TextCell 50 2 "Maker"
Set-TextCell 50 3 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell 50 4 "1.301.46"
Set-TextCell 50 5 "  -2.85%  "
Set-TextCell 51 2 "MXToken"
Set-TextCell 51 3 "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell 51 4 "2.91"
Set-TextCell 51 5 "  +0.89%  "
